# The field " m:forVar.name " (a Word field made of begin/instrText/end
# runs) is rewritten as plain literal text runs:
#   "{" "m" ":" "for" "Var" ".name}"
# keeping the orange theme color (accent6, 50% darker) on "for" and "Var".

$d = $word.ActiveDocument

# This template contains exactly one field: the " m:forVar.name " field.
$field = $d.Fields.Item(1)

# Remember where the field starts (one character before its instruction
# text starts is the position of the "begin" field-char run) so we can
# insert the replacement text at the same spot once the field is gone.
$insertAt = $field.Code.Start - 1

# Removes all of the field's runs (fldChar begin/end and instrText runs).
$field.Delete()

$target = $d.Range($insertAt, $insertAt)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' +
          '<w:r><w:t>{</w:t></w:r>' +
          '<w:r><w:t>m</w:t></w:r>' +
          '<w:r><w:t>:</w:t></w:r>' +
          '<w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>for</w:t></w:r>' +
          '<w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>Var</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve">.name}</w:t></w:r>' +
        '</w:p></w:body>' +
      '</w:document>' +
    '</pkg:xmlData></pkg:part>' +
  '</pkg:package>'

$target.InsertXML($xml)
